$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '37.518.84'
$ws.Range('E2').Value = '  +1.52%  '

# Row 3
$ws.Range('D3').Value = '2.032.55'
$ws.Range('E3').Value = '  +2.38%  '

# Row 4
$ws.Range('E4').Value = '  -0.11%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '256.25'
$ws.Range('E5').Value = '  +4.45%  '

# Row 6
$ws.Range('E6').Value = '  -0.74%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.08%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '57.42'
$ws.Range('E8').Value = '  -5.88%  '

# Row 9
$ws.Range('E9').Value = '  +1.01%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0795'
$ws.Range('E10').Value = '  -0.57%  '

# Row 11
$ws.Range('E11').Value = '  -1.55%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.77'
$ws.Range('E12').Value = '  -1.54%  '

# Row 13
$ws.Range('D13').Value = '2.331.76'
$ws.Range('E13').Value = '  +2.46%  '

# Row 14
$ws.Range('E14').Value = '  -2.91%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.39'
$ws.Range('E15').Value = '  -3.61%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.37'
$ws.Range('E16').Value = '  -2.01%  '

# Row 17
$ws.Range('D17').Value = '2.034.97'
$ws.Range('E17').Value = '  +2.43%  '

# Row 18
$ws.Range('D18').Value = '37.422.36'
$ws.Range('E18').Value = '  +1.49%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.98'
$ws.Range('E19').Value = '  -0.45%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0855'
$ws.Range('E20').Value = '  -0.95%  '

# Row 21
$ws.Range('E21').Value = '  +0.49%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '229.25'
$ws.Range('E22').Value = '  -0.37%  '

# Row 23
$ws.Range('E23').Value = '  +4.92%  '

# Row 24
$ws.Range('E24').Value = '  -0.01%  '

# Row 25
$ws.Range('E25').Value = '  -0.98%  '

# Row 26
$ws.Range('E26').Value = '  -1.72%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '163.43'
$ws.Range('E27').Value = '  +0.13%  '

# Row 28
$ws.Range('E28').Value = '  -7.90%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.96'
$ws.Range('E29').Value = '  +1.98%  '

# Row 30
$ws.Range('E30').Value = '  -1.16%  '

# Row 31
$ws.Range('E31').Value = '  -1.13%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0667'
$ws.Range('E32').Value = '  +7.28%  '

# Row 33
$ws.Range('E33').Value = '  -3.46%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.56'
$ws.Range('E34').Value = '  +0.22%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.46'
$ws.Range('E35').Value = '  +7.59%  '

# Row 36
$ws.Range('E36').Value = '  +0.02%  '

# Row 37
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.82'
$ws.Range('E37').Value = '  +1.97%  '

# Row 38
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.40'
$ws.Range('E38').Value = '  +1.13%  '

# Row 39
$ws.Range('E39').Value = '  -3.52%  '

# Row 40
$ws.Range('E40').Value = '  +3.83%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0966'
$ws.Range('E41').Value = '  -3.11%  '

# Row 42
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.19'
$ws.Range('E42').Value = '  +1.00%  '

# Row 43
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0216'
$ws.Range('E43').Value = '  +1.26%  '

# Row 44
$ws.Range('D44').Value = '1.401.91'
$ws.Range('E44').Value = '  +2.18%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '16.04'
$ws.Range('E45').Value = '  -3.47%  '

# Row 46
$ws.Range('E46').Value = '  +0.57%  '

# Row 47
$ws.Range('E47').Value = '  +0.77%  '

# Row 48
$ws.Range('E48').Value = '  +0.89%  '

# Row 49
$ws.Range('E49').Value = '  +1.77%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.02'
$ws.Range('E50').Value = '  +1.67%  '

# Row 51
$ws.Range('D51').Value = '2.222.58'
$ws.Range('E51').Value = '  +2.45%  '
